# ============================================================================
# Updates "杭州-漫展信息.xlsx" to the next scrape snapshot (gh-pages output
# generated at 456a3b4):
#   - Sheet "展览"   : insert 2 new events, refresh "想去人数" (views) counts
#   - Sheet "演出"   : insert 1 new event, refresh a views count
#   - Sheet "全部类型": refresh views counts only (no new rows here yet)
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet 1: 展览 (exhibitions)
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# --- insert the new "白日梦次元动漫嘉年华" event before the old row 17 ------
$ws1.Rows.Item(17).Insert()
$ws1.Cells.Item(17,2).Value2 = "2024-06-23"
$ws1.Cells.Item(17,3).Value2 = "杭州·第二届白日梦次元动漫嘉年华"
$ws1.Cells.Item(17,4).Value2 = "康候圣街99号 顺丰创新中心"
$ws1.Cells.Item(17,5).Value2 = "2024.06.23 10:00-06.23 17:00"
$ws1.Cells.Item(17,6).Value2 = 0
$ws1.Cells.Item(17,7).Value2 = 68
$ws1.Cells.Item(17,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=86307"
$ws1.Cells.Item(17,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202405/qHcyIUL31715752173541.jpeg"

# --- insert the new "第四届ArknightsOnly" event (ends up at row 33) --------
$ws1.Rows.Item(33).Insert()
$ws1.Cells.Item(33,2).Value2 = "2024-07-20"
$ws1.Cells.Item(33,3).Value2 = "杭州·第四届ArknightsOnly·狼与黑荆棘（明日方舟Only）"
$ws1.Cells.Item(33,4).Value2 = "康候圣街99号 顺丰创新中心"
$ws1.Cells.Item(33,5).Value2 = "2024.07.20 10:00-07.20 17:00"
$ws1.Cells.Item(33,6).Value2 = 0
$ws1.Cells.Item(33,7).Value2 = 79
$ws1.Cells.Item(33,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=86305"
$ws1.Cells.Item(33,9).Value2 = "//i1.hdslb.com/bfs/openplatform/202405/cpoiCink1716554216810.png"

# --- fix up formatting of the two newly-inserted index cells (column A) ----
foreach ($r in @(17,33)) {
    $cell = $ws1.Cells.Item($r,1)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(1).LineStyle = 1
    $cell.Borders.Item(2).LineStyle = 1
    $cell.Borders.Item(3).LineStyle = 1
    $cell.Borders.Item(4).LineStyle = 1
}

# --- refreshed "想去人数" (views) counts, addressed by their FINAL row ------
$sheet1Views = @{
    2  = 568
    3  = 5342
    8  = 369
    12 = 3057
    13 = 1897
    16 = 186
    18 = 131
    20 = 966
    21 = 346
    22 = 45
    23 = 3499
    24 = 1098
    25 = 2784
    27 = 1752
    28 = 4008
    31 = 457
    36 = 1255
    37 = 56
    38 = 1015
    39 = 660
    40 = 504
    41 = 402
    42 = 309
}
foreach ($row in $sheet1Views.Keys) {
    $ws1.Cells.Item($row,6).Value2 = $sheet1Views[$row]
}

# --- renumber the "#" index column (A) now that rows shifted ---------------
for ($r = 2; $r -le 43; $r++) {
    $ws1.Cells.Item($r,1).Value2 = $r - 1
}

# ----------------------------------------------------------------------------
# Sheet 2: 演出 (performances)
# ----------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# --- insert the new "Marcin Patrzalek" concert event ------------------------
$ws2.Rows.Item(24).Insert()
$ws2.Cells.Item(24,2).Value2 = "2024-08-16"
$ws2.Cells.Item(24,3).Value2 = "杭州·Marcin Patrzalek 2024 《原声之龙》指弹吉他音乐会"
$ws2.Cells.Item(24,4).Value2 = "新北街85号三层 66livehouse"
$ws2.Cells.Item(24,5).Value2 = "2024.08.16 19:30-08.16 21:00"
$ws2.Cells.Item(24,6).Value2 = 0
$ws2.Cells.Item(24,7).Value2 = 380
$ws2.Cells.Item(24,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=86294"
$ws2.Cells.Item(24,9).Value2 = "//i2.hdslb.com/bfs/openplatform/202405/sbxCQPzr1716547720261.jpeg"

$cell2 = $ws2.Cells.Item(24,1)
$cell2.Font.Bold = $true
$cell2.HorizontalAlignment = -4108
$cell2.VerticalAlignment = -4160
$cell2.Borders.Item(1).LineStyle = 1
$cell2.Borders.Item(2).LineStyle = 1
$cell2.Borders.Item(3).LineStyle = 1
$cell2.Borders.Item(4).LineStyle = 1

# --- refreshed views count for the shifted "动漫作品主题音乐会" row --------
$ws2.Cells.Item(25,6).Value2 = 5

# --- renumber the "#" index column (A) ---------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws2.Cells.Item($r,1).Value2 = $r - 1
}

# ----------------------------------------------------------------------------
# Sheet 4: 全部类型 (all types) - refresh views counts only, rows unchanged
# ----------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Views = @{
    2  = 568
    3  = 568
    4  = 5342
    9  = 369
    11 = 3057
    13 = 1897
    17 = 186
    21 = 131
    22 = 966
    23 = 346
    24 = 3499
    27 = 1098
    28 = 2784
    29 = 1752
    30 = 4008
    38 = 1255
    39 = 56
    40 = 1015
    42 = 660
    44 = 402
    47 = 5
    48 = 309
}
foreach ($row in $sheet4Views.Keys) {
    $ws4.Cells.Item($row,6).Value2 = $sheet4Views[$row]
}
